# columns_number_parameters.xlsx (HU) update
# - B27 (row for "columnsPartnershipU1b") changes from the number 27 to the
#   text value "25" (quote-prefixed text, matching the new shared string
#   "25" and the new quotePrefix cell style).
# - The active selection moves from B41 to C27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# Leading apostrophe forces Excel to store this as quote-prefixed text
# rather than a number, matching the target workbook's new shared string
# "25" plus the new cellXfs entry with quotePrefix="1".
$ws.Range("B27").Value = "'25"

# Update the sheet's active cell/selection to C27.
$ws.Range("C27").Select()
